$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D contains price strings that often look numeric (e.g. "5.25").
# Force the whole data range to Text format first so Excel stores the
# values as literal strings instead of converting/rounding them as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.860.39"
$ws.Range("E2").Value = "  +1.70%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.622.52"
$ws.Range("E3").Value = "  +1.75%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.08%  "

# Row 5 - BNB
$ws.Range("D5").Value = "605.49"
$ws.Range("E5").Value = "  +2.33%  "

# Row 6 - Solana
$ws.Range("D6").Value = "154.58"
$ws.Range("E6").Value = "  +0.54%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +2.25%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.621.72"
$ws.Range("E9").Value = "  +1.81%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.129"
$ws.Range("E10").Value = "  +14.10%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.85%  "

# Row 12 - Toncoin
$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  +1.59%  "

# Row 13 - Cardano
$ws.Range("D13").Value = "0.355"
$ws.Range("E13").Value = "  +0.32%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "27.83"
$ws.Range("E14").Value = "  -0.45%  "

# Row 15 - ShibaInu
$ws.Range("D15").Value = "0.0000189"
$ws.Range("E15").Value = "  +5.52%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "3.107.20"
$ws.Range("E16").Value = "  +1.88%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "67.815.12"
$ws.Range("E17").Value = "  +1.63%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.615.89"
$ws.Range("E18").Value = "  +1.15%  "

# Row 19 - was BitcoinCash, now Chainlink
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "11.24"
$ws.Range("E19").Value = "  +0.47%  "

# Row 20 - was Chainlink, now BitcoinCash
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "366.83"
$ws.Range("E20").Value = "  +3.89%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "7.68"
$ws.Range("E21").Value = "  -0.52%  "

# Row 22 - Polkadot
$ws.Range("E22").Value = "  -0.20%  "

# Row 23 - SuiNetwork
$ws.Range("D23").Value = "2.06"
$ws.Range("E23").Value = "  +2.27%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.09%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "70.42"
$ws.Range("E25").Value = "  +4.75%  "

# Row 26 - Aptos
$ws.Range("D26").Value = "9.92"
$ws.Range("E26").Value = "  -3.22%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  +3.27%  "

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "2.746.35"
$ws.Range("E28").Value = "  +1.39%  "

# Row 29 - Bittensor
$ws.Range("D29").Value = "580.71"
$ws.Range("E29").Value = "  -1.95%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  +0.04%  "

# Row 31 - Fetch.AI
$ws.Range("D31").Value = "1.44"
$ws.Range("E31").Value = "  -0.53%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "7.93"
$ws.Range("E32").Value = "  -0.30%  "

# Row 33 - PancakeSwap
$ws.Range("D33").Value = "1.88"
$ws.Range("E33").Value = "  +1.82%  "

# Row 34 - Kaspa
$ws.Range("E34").Value = "  -0.76%  "

# Row 35 - FirstDigitalUSD
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.03%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  -1.75%  "

# Row 37 - NEARProtocol
$ws.Range("D37").Value = "4.97"
$ws.Range("E37").Value = "  +0.31%  "

# Row 38 - Monero
$ws.Range("D38").Value = "157.63"
$ws.Range("E38").Value = "  +2.69%  "

# Row 39 - EthereumClassic
$ws.Range("D39").Value = "19.46"
$ws.Range("E39").Value = "  +2.05%  "

# Row 40 - PolygonEcosystemToken
$ws.Range("D40").Value = "0.371"
$ws.Range("E40").Value = "  +1.01%  "

# Row 41 - RenderToken
$ws.Range("D41").Value = "5.38"
$ws.Range("E41").Value = "  -0.53%  "

# Row 42 - Stacks
$ws.Range("D42").Value = "1.86"
$ws.Range("E42").Value = "  +4.53%  "

# Row 43 - dogwifhat
$ws.Range("D43").Value = "2.64"
$ws.Range("E43").Value = "  +1.62%  "

# Row 44 - OKB
$ws.Range("D44").Value = "41.17"
$ws.Range("E44").Value = "  -0.63%  "

# Row 45 - USDe
$ws.Range("E45").Value = "  +0.04%  "

# Row 46 - WhiteBITCoin
$ws.Range("D46").Value = "16.44"
$ws.Range("E46").Value = "  +0.21%  "

# Row 47 - Aave
$ws.Range("D47").Value = "157.39"
$ws.Range("E47").Value = "  +1.79%  "

# Row 48 - BabyDogeCoin
$ws.Range("D48").Value = "0.0₆0290"
$ws.Range("E48").Value = "  -5.07%  "

# Row 49 - Filecoin
$ws.Range("E49").Value = "  +1.00%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "21.04"
$ws.Range("E50").Value = "  +0.01%  "

# Row 51 - was Hedera, now Mantle
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.625"
$ws.Range("E51").Value = "  +2.38%  "

# Restore the default (no explicit) style on column D now that the text
# values are safely stored, matching the original workbook's formatting.
$ws.Range("D2:D51").Style = "Normal"
